# Lecture 16 / "Simple Buffer" slide: drop the "Buffer Full?" / "Buffer
# Empty?" bullets and split the "Concept of circular buffer" bullet into
# two runs ("Concept " + "of circular buffer").

$p = $ppt.ActivePresentation

# Locate the shape that holds the "Buffer Full?" bullet text (slide 5,
# "Content Placeholder 36") without hard-coding slide/shape numbers.
$targetShape = $null
foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText -and ($sh.TextFrame.TextRange.Text -like "*Buffer Full?*")) {
                $targetShape = $sh
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Remove the "Buffer Full?" and "Buffer Empty?" paragraphs entirely. Walk
# backwards so deleting a paragraph doesn't shift the index of paragraphs
# that still need to be visited. (TextRange.Text on a paragraph carries a
# trailing CR, hence the TrimEnd().)
for ($i = $tr.Paragraphs().Count; $i -ge 1; $i--) {
    $txt = $tr.Paragraphs($i).Text.TrimEnd()
    if ($txt -eq "Buffer Full?" -or $txt -eq "Buffer Empty?") {
        $tr.Paragraphs($i).Delete()
    }
}

# Split "Concept of circular buffer" into two runs: "Concept " and
# "of circular buffer" (same formatting, just two separate <a:r> runs).
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text.TrimEnd() -eq "Concept of circular buffer") {
        $firstRun = $para.Characters(1, 8)
        # Touching the run's formatting forces PowerPoint to split the
        # paragraph's single run into two separate runs at this boundary.
        $firstRun.Font.Bold = $firstRun.Font.Bold
    }
}

Write-Output $tr.Text
